$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3495
$ws1.Range("F6").Value = 426
$ws1.Range("F9").Value = 51
$ws1.Range("F12").Value = 1722
$ws1.Range("F13").Value = 128

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3495
$ws4.Range("F6").Value = 426
$ws4.Range("F10").Value = 51
$ws4.Range("F15").Value = 1722
$ws4.Range("F16").Value = 128
